$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(902, 1, 0, 0, 0, 0),
  @(401, 9, 48, 67, 75, 45),
  @(901, 16, 15, 45, 60, 60),
  @(501, 9, 52, 30, 75, 45),
  @(701, 3, 90, 45, 97, 15),
  @(1202, 2, 10, 10, 10, 10),
  @(301, 6, 45, 30, 60, 45),
  @(201, 9, 30, 15, 45, 30),
  @(801, 3, 67, 65, 52, 45),
  @(101, 9, 30, 15, 60, 15),
  @(1001, 18, 30, 75, 60, 72),
  @(601, 9, 60, 67, 60, 42),
  @(1201, 2, 10, 10, 10, 10),
  @(1203, 3, 15, 15, 15, 15),
  @(1101, 0, 15, 30, 30, 0),
  @(1, 0, 2, 2, 2, 2),
  @(3, 0, 3, 3, 3, 3),
  @(802, 0, 4, 5, 4, 0),
  @(502, 0, 4, 0, 0, 0),
  @(2, 0, 2, 2, 2, 2)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
  $row = $startRow + $i
  $vals = $data[$i]
  for ($c = 0; $c -lt $vals.Length; $c++) {
    $ws.Cells.Item($row, $c + 1).Value = $vals[$c]
  }
}
